# Updated cryptos list on Mon Aug 28 20:22:32 UTC 2023 with GitHub Actions
# Refreshes price/volume figures (Price + Volume(1h) columns) on the active sheet,
# including a few coins whose rank/row order changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.104.98'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.652.76'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.72'
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5251'
$ws.Range("E6").Value = '  -0.94%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2664'
$ws.Range("E8").Value = '  +0.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06355'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.56'
$ws.Range("E10").Value = '  -1.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07692'
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.771.25'
$ws.Range("E12").Value = '  +5.59%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.598'
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.879.48'
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5602'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8218'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.37'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.103.08'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.695'
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '191.03'
$ws.Range("E22").Value = '  -3.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.978'
$ws.Range("E23").Value = '  -1.28%  '
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.75'
$ws.Range("E25").Value = '  -0.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1199'
$ws.Range("E26").Value = '  -1.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.256'
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.92'
$ws.Range("E28").Value = '  -1.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.495'
$ws.Range("E29").Value = '  -1.19%  '
$ws.Range("E30").Value = '  -4.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.271'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.497'
$ws.Range("E32").Value = '  -1.07%  '
$ws.Range("E33").Value = '  +1.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.580'
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.793'
$ws.Range("E35").Value = '  -1.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9471'
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.409'
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5778'
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.972'
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8401'
$ws.Range("E42").Value = '  -2.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.019.71'
$ws.Range("E43").Value = '  -5.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.61'
$ws.Range("E44").Value = '  -1.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.791.06'
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.29'
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈106'
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("E48").Value = '  -1.12%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05328'
$ws.Range("E49").Value = '  +3.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.017'
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4342'
$ws.Range("E51").Value = '  -1.59%  '

